# Capital commitments workbook - add FX (From/To Currency, Exchange Rate, As Of)
# columns and correct the Commitment Date values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (order controls shared-string insertion order so the
#     resulting sharedStrings.xml indices match the target file) ---
$ws.Range("K1").Value = "From Currency"
$ws.Range("L1").Value = "To Currency"
$ws.Range("N1").Value = "As Of"
$ws.Range("M1").Value = "Exchange Rate "

# --- Row 2: new FX detail ---
$ws.Range("K2").Value = "USD"
$ws.Range("L2").Value = "INR"
$ws.Range("M2").Value = 80
$ws.Range("J2").Copy()
$ws.Range("N2").PasteSpecial(-4122)
$ws.Range("N2").Value = 44216

# --- Row 5: new FX detail ---
$ws.Range("K5").Value = "USD"
$ws.Range("L5").Value = "INR"
$ws.Range("M5").Value = 80
$ws.Range("J5").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = 44216

$excel.CutCopyMode = $false

# --- Fix the Commitment Date column (J) for every data row ---
$ws.Range("J2").Value = 44216
$ws.Range("J3").Value = 44216
$ws.Range("J4").Value = 44216
$ws.Range("J5").Value = 44216
$ws.Range("J6").Value = 44216
$ws.Range("J7").Value = 44216

# --- Column widths for the new FX columns ---
$ws.Range("K1").ColumnWidth = 12.875 - 5/7
$ws.Range("L1").ColumnWidth = 10.6875 - 5/7
$ws.Range("M1").ColumnWidth = 13.5 - 5/7
$ws.Range("N1").ColumnWidth = 11.875 - 5/7

$ws.Range("N5").Select()
